$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46001
$ws.Range("B2").Value = 46001.77383101852
$ws.Range("C2").Value = 'Stock Market Today: Dow Rises As Fed, Powell Loom; Palantir Gains (Live Coverage)'
$ws.Range("D2").Value = 'https://finance.yahoo.com/m/d2a83895-3dec-3c32-978d-30d95c7a44de/stock-market-today%3A-dow-rises.html'
$ws.Range("E2").Value = 'Investor''s Business Daily'
$ws.Range("F2").Value = 'Yahoo_Search'

$ws.Range("A3").Value = 46001
$ws.Range("B3").Value = 46001.77083333334
$ws.Range("C3").Value = 'Private markets: Thinking outside the traditional 60/40 portfolio'
$ws.Range("D3").Value = 'https://finance.yahoo.com/video/private-markets-thinking-outside-traditional-183000929.html'
$ws.Range("E3").Value = 'Yahoo Finance Video'
$ws.Range("F3").Value = 'Yahoo_Search'

$ws.Range("A4").Value = 46001
$ws.Range("B4").Value = 46001.76931712963
$ws.Range("C4").Value = 'Stock market today: Dow, S&P 500, Nasdaq diverge with Fed rate decision on deck'
$ws.Range("D4").Value = 'https://finance.yahoo.com/news/live/stock-market-today-dow-sp-500-nasdaq-diverge-with-fed-rate-decision-on-deck-235110803.html'
$ws.Range("E4").Value = 'Yahoo Finance'
$ws.Range("F4").Value = 'Yahoo_Search'

$ws.Range("A5").Value = 46001
$ws.Range("B5").Value = 46001.75576388889
$ws.Range("C5").Value = 'Exchange-Traded Funds Climb as US Equities Mixed After Midday'
$ws.Range("D5").Value = 'https://finance.yahoo.com/news/exchange-traded-funds-climb-us-180818818.html'
$ws.Range("E5").Value = 'MT Newswires'
$ws.Range("F5").Value = 'Yahoo_Search'

$ws.Range("A6").Value = 46001
$ws.Range("B6").Value = 46001.73265046296
$ws.Range("C6").Value = 'US Equity Indexes Trade Mixed Ahead of Anticipated Divisions in Fed Policy'
$ws.Range("D6").Value = 'https://finance.yahoo.com/news/us-equity-indexes-trade-mixed-173501076.html'
$ws.Range("E6").Value = 'MT Newswires'
$ws.Range("F6").Value = 'Yahoo_Search'

$ws.Range("A7").Value = 46001
$ws.Range("B7").Value = 46001.71746527778
$ws.Range("C7").Value = 'Stocks mark time ahead of Fed decision'
$ws.Range("D7").Value = 'https://finance.yahoo.com/news/asian-stocks-retreat-traders-eye-025020585.html'
$ws.Range("E7").Value = 'AFP'
$ws.Range("F7").Value = 'Yahoo_Search'

$ws.Range("A8").Value = 46001
$ws.Range("B8").Value = 46001.71572916667
$ws.Range("C8").Value = 'Talk of an AI bubble is just ''ridiculous,'' this strategist says'
$ws.Range("D8").Value = 'https://finance.yahoo.com/video/talk-ai-bubble-just-ridiculous-171039241.html'
$ws.Range("E8").Value = 'Yahoo Finance Video'
$ws.Range("F8").Value = 'Yahoo_Search'

$ws.Range("A9").Value = 46001
$ws.Range("B9").Value = 46001.6916087963
$ws.Range("C9").Value = 'Top Midday Stories: DeepSeek Said to Be Using Smuggled Nvidia Chips; Meta Pivoting Toward Money-Making AI Model Over Open Source'
$ws.Range("D9").Value = 'https://finance.yahoo.com/news/top-midday-stories-deepseek-said-163555839.html'
$ws.Range("E9").Value = 'MT Newswires'
$ws.Range("F9").Value = 'Yahoo_Search'

$ws.Range("A10").Value = 46001
$ws.Range("B10").Value = 46001.66736111111
$ws.Range("C10").Value = 'Markets Brace for a Hawkish Fed'
$ws.Range("D10").Value = 'https://finance.yahoo.com/news/markets-brace-hawkish-fed-160100708.html'
$ws.Range("E10").Value = 'Zacks'
$ws.Range("F10").Value = 'Yahoo_Search'

$ws.Range("A11").Value = 46001
$ws.Range("B11").Value = 46001.64652777778
$ws.Range("C11").Value = 'Pre-markets Remain "Wait & See" Ahead of Fed Decision'
$ws.Range("D11").Value = 'https://finance.yahoo.com/news/pre-markets-remain-wait-see-153100353.html'
$ws.Range("E11").Value = 'Zacks'
$ws.Range("F11").Value = 'Yahoo_Search'

$ws.Range("A22").Value = 46001
$ws.Range("B22").Value = 46001.77325231482
$ws.Range("C22").Value = 'Trump says Nvidia can sell H200 chips to China, but nothing is guaranteed'
$ws.Range("D22").Value = 'https://finance.yahoo.com/news/trump-says-nvidia-can-sell-h200-chips-to-china-but-nothing-is-guaranteed-183329002.html'
$ws.Range("E22").Value = 'Yahoo Finance'
$ws.Range("F22").Value = 'Yahoo_Search'

$ws.Range("A23").Value = 46001
$ws.Range("B23").Value = 46001.77228009259
$ws.Range("C23").Value = 'BioNTech and BMS’ bispecific touts 59.3% nine-month PFS in TNBC'
$ws.Range("D23").Value = 'https://finance.yahoo.com/m/aa32fe1c-edb1-320f-833b-809a5f77f186/biontech-and-bms%E2%80%99-bispecific.html'
$ws.Range("E23").Value = 'Clinical Trials Arena'
$ws.Range("F23").Value = 'Yahoo_Search'

$ws.Range("A24").Value = 46001
$ws.Range("B24").Value = 46001.77152777778
$ws.Range("C24").Value = 'Totem Launches WeFunder Community Round to Scale Offline, Decentralized Friend-Finding Technology'
$ws.Range("D24").Value = 'https://finance.yahoo.com/news/totem-launches-wefunder-community-round-183100370.html'
$ws.Range("E24").Value = 'Business Wire'
$ws.Range("F24").Value = 'Yahoo_Search'

$ws.Range("A25").Value = 46001
$ws.Range("B25").Value = 46001.77084490741
$ws.Range("C25").Value = 'Betting markets predict Time''s 2025 Person of the Year won''t be a person at all'
$ws.Range("D25").Value = 'https://finance.yahoo.com/m/1c14bd52-3022-3610-b772-b21c516c708e/betting-markets-predict.html'
$ws.Range("E25").Value = 'Business Insider'
$ws.Range("F25").Value = 'Yahoo_Search'

$ws.Range("A26").Value = 46001
$ws.Range("B26").Value = 46001.77083333334
$ws.Range("C26").Value = 'New Podcast Discusses Strengthening Military Families With Innovation and Partnership'
$ws.Range("D26").Value = 'https://finance.yahoo.com/news/podcast-discusses-strengthening-military-families-183000468.html'
$ws.Range("E26").Value = 'ACCESS Newswire'
$ws.Range("F26").Value = 'Yahoo_Search'

$ws.Range("A27").Value = 46001
$ws.Range("B27").Value = 46001.77083333334
$ws.Range("C27").Value = 'A2Z Cust2Mate to Host Investor Virtual Meeting on Tuesday, December 16, 2025'
$ws.Range("D27").Value = 'https://finance.yahoo.com/news/a2z-cust2mate-host-investor-virtual-183000417.html'
$ws.Range("E27").Value = 'ACCESS Newswire'
$ws.Range("F27").Value = 'Yahoo_Search'

$ws.Range("A28").Value = 46001
$ws.Range("B28").Value = 46001.77083333334
$ws.Range("C28").Value = 'Altria Declares Regular Quarterly Dividend of $1.06 Per Share'
$ws.Range("D28").Value = 'https://finance.yahoo.com/news/altria-declares-regular-quarterly-dividend-183000365.html'
$ws.Range("E28").Value = 'Business Wire'
$ws.Range("F28").Value = 'Yahoo_Search'

$ws.Range("A29").Value = 46001
$ws.Range("B29").Value = 46001.77083333334
$ws.Range("C29").Value = 'What 10 Years and 600 Studies Tell Us About Sustainable Business Success'
$ws.Range("D29").Value = 'https://finance.yahoo.com/news/10-years-600-studies-tell-183000749.html'
$ws.Range("E29").Value = 'ACCESS Newswire'
$ws.Range("F29").Value = 'Yahoo_Search'

$ws.Range("A30").Value = 46001
$ws.Range("B30").Value = 46001.77083333334
$ws.Range("C30").Value = 'Carrick Therapeutics Announces Positive Results from Phase 2 Randomized Trial of Samuraciclib in Combination with Fulvestrant in Patients with Hormone Receptor Positive, HER2 Negative Advanced Breast Cancer'
$ws.Range("D30").Value = 'https://finance.yahoo.com/news/carrick-therapeutics-announces-positive-results-183000909.html'
$ws.Range("E30").Value = 'GlobeNewswire'
$ws.Range("F30").Value = 'Yahoo_Search'

$ws.Range("A31").Value = 45995
$ws.Range("B31").Value = 45995.67222222222
$ws.Range("C31").Value = 'ETFs to Keep Your Portfolio on Track in the Long Term'
$ws.Range("D31").Value = 'https://finance.yahoo.com/news/etfs-keep-portfolio-track-long-160800229.html'
$ws.Range("E31").Value = 'Zacks'
$ws.Range("F31").Value = 'Yahoo_Search'

$ws.Range("A32").Value = 45995
$ws.Range("B32").Value = 45995.65972222222
$ws.Range("C32").Value = 'Here''s Why This ETF Is a Multimillionaire Maker'
$ws.Range("D32").Value = 'https://finance.yahoo.com/m/41bb49e9-5c08-38ac-8d89-c0d326db1a26/here%27s-why-this-etf-is-a.html'
$ws.Range("E32").Value = 'Motley Fool'
$ws.Range("F32").Value = 'Yahoo_Search'

$ws.Range("A33").Value = 45985
$ws.Range("B33").Value = 45985.34375
$ws.Range("C33").Value = 'If I Could Choose Only 1 ETF to Buy and Hold Forever, This Would Be It'
$ws.Range("D33").Value = 'https://finance.yahoo.com/m/d1fdddb9-c08d-33c4-a6a5-23354a1e8a76/if-i-could-choose-only-1-etf.html'
$ws.Range("E33").Value = 'Motley Fool'
$ws.Range("F33").Value = 'Yahoo_Search'

$ws.Range("A34").Value = 45984
$ws.Range("B34").Value = 45984.86805555555
$ws.Range("C34").Value = 'Here''s How Many Shares of the Vanguard Total Stock Market ETF (VTI) You''d Need for $500 in Yearly Dividends'
$ws.Range("D34").Value = 'https://finance.yahoo.com/m/34c30e63-4b19-3dde-83ca-9034e11f82ce/here%27s-how-many-shares-of-the.html'
$ws.Range("E34").Value = 'Motley Fool'
$ws.Range("F34").Value = 'Yahoo_Search'

$ws.Range("A35").Value = 45983
$ws.Range("B35").Value = 45983.18337962963
$ws.Range("C35").Value = 'Buy These 4 ETFs if You Want to be Rich in 2026, According to John Liang'
$ws.Range("D35").Value = 'https://finance.yahoo.com/news/buy-4-etfs-want-rich-042404092.html'
$ws.Range("E35").Value = 'GOBankingRates'
$ws.Range("F35").Value = 'Yahoo_Search'

$ws.Range("A36").Value = 45980
$ws.Range("B36").Value = 45980.98965277777
$ws.Range("C36").Value = 'Two Paths to the Total US Market: Vanguard’s Total Stock Market ETF vs. iShares’ Core S&P Total US Stock Market ETF'
$ws.Range("D36").Value = 'https://finance.yahoo.com/m/7dd78558-e7d1-32b6-95bf-8247fbbeb672/two-paths-to-the-total-us.html'
$ws.Range("E36").Value = 'Motley Fool'
$ws.Range("F36").Value = 'Yahoo_Search'

$ws.Range("A37").Value = 45980
$ws.Range("B37").Value = 45980.875625
$ws.Range("C37").Value = 'The Clear Winner for Building Long Term Wealth, QQQ or VTI?'
$ws.Range("D37").Value = 'https://finance.yahoo.com/m/794b4311-7b64-3873-aae3-e5bc37477c25/the-clear-winner-for-building.html'
$ws.Range("E37").Value = '24/7 Wall St.'
$ws.Range("F37").Value = 'Yahoo_Search'

$ws.Range("A38").Value = 45978
$ws.Range("B38").Value = 45978.46458333333
$ws.Range("C38").Value = 'This 16% Difference Could Make the Vanguard Total Stock Market ETF Outperform the S&P 500 During a Stock Market Sell-Off'
$ws.Range("D38").Value = 'https://finance.yahoo.com/m/ab16c17e-e4d2-381d-8e3d-987a9323e2f0/this-16%25-difference-could.html'
$ws.Range("E38").Value = 'Motley Fool'
$ws.Range("F38").Value = 'Yahoo_Search'

$ws.Range("A39").Value = 45975
$ws.Range("B39").Value = 45975.65974537037
$ws.Range("C39").Value = 'Stop Panicking and Look at the Charts: What Barchart’s Technical Indicators Are Telling Us Now'
$ws.Range("D39").Value = 'https://finance.yahoo.com/m/df3960f4-aac8-328e-b3fa-0e8a83a609ed/stop-panicking-and-look-at.html'
$ws.Range("E39").Value = 'Barchart'
$ws.Range("F39").Value = 'Yahoo_Search'

$ws.Range("A40").Value = 45960
$ws.Range("B40").Value = 45960.36111111111
$ws.Range("C40").Value = 'The Best Vanguard ETF to Invest $2,000 in Right Now'
$ws.Range("D40").Value = 'https://finance.yahoo.com/m/3925f7af-f7bb-37c6-825f-3b46e1186c96/the-best-vanguard-etf-to.html'
$ws.Range("E40").Value = 'Motley Fool'
$ws.Range("F40").Value = 'Yahoo_Search'

$ws.Range("A41").Value = 46001
$ws.Range("B41").Value = 46001.77708333333
$ws.Range("C41").Value = 'Photonic Sets New Standard with Distributed Quantum Resource Estimation'
$ws.Range("D41").Value = 'https://finance.yahoo.com/news/photonic-sets-standard-distributed-quantum-183900917.html'
$ws.Range("E41").Value = 'GlobeNewswire'
$ws.Range("F41").Value = 'Yahoo_Search'

$ws.Range("A42").Value = 46001
$ws.Range("B42").Value = 46001.7769212963
$ws.Range("C42").Value = 'Athene exec reveals the 2 mistakes retirees make — and the kicker that could eviscerate savings'
$ws.Range("D42").Value = 'https://finance.yahoo.com/news/athene-exec-reveals-the-2-mistakes-retirees-make--and-the-kicker-that-could-eviscerate-savings-183846206.html'
$ws.Range("E42").Value = 'Yahoo Finance'
$ws.Range("F42").Value = 'Yahoo_Search'

$ws.Range("A43").Value = 46001
$ws.Range("B43").Value = 46001.77638888889
$ws.Range("C43").Value = '31 Japanese Startups to Exhibit at CES 2026 Japan Pavilion Including 4 Innovation Award Winners'
$ws.Range("D43").Value = 'https://finance.yahoo.com/news/31-japanese-startups-exhibit-ces-183800738.html'
$ws.Range("E43").Value = 'Business Wire'
$ws.Range("F43").Value = 'Yahoo_Search'

$ws.Range("A44").Value = 46001
$ws.Range("B44").Value = 46001.77569444444
$ws.Range("C44").Value = '2 Cruise Line Stocks Are Moving in Different Directions'
$ws.Range("D44").Value = 'https://finance.yahoo.com/m/654ead26-dbc7-307e-8491-137b3cf276c3/2-cruise-line-stocks-are.html'
$ws.Range("E44").Value = 'Motley Fool'
$ws.Range("F44").Value = 'Yahoo_Search'

$ws.Range("A45").Value = 46001
$ws.Range("B45").Value = 46001.775
$ws.Range("C45").Value = 'AltaSea Taps LA City & LAEDC Veteran to Join Leadership Team'
$ws.Range("D45").Value = 'https://finance.yahoo.com/news/altasea-taps-la-city-laedc-183600829.html'
$ws.Range("E45").Value = 'Business Wire'
$ws.Range("F45").Value = 'Yahoo_Search'

$ws.Range("A46").Value = 46001
$ws.Range("B46").Value = 46001.775
$ws.Range("C46").Value = 'Tickets Now on Sale for 2026 International Fireworks Championship in Las Vegas'
$ws.Range("D46").Value = 'https://finance.yahoo.com/news/tickets-now-sale-2026-international-183600777.html'
$ws.Range("E46").Value = 'Business Wire'
$ws.Range("F46").Value = 'Yahoo_Search'

$ws.Range("A47").Value = 46001
$ws.Range("B47").Value = 46001.77494212963
$ws.Range("C47").Value = 'These Are The 5 Best Stocks To Buy Now Or Watch'
$ws.Range("D47").Value = 'https://finance.yahoo.com/m/5f695c14-bc91-363c-995e-e994c1f0807e/these-are-the-5-best-stocks.html'
$ws.Range("E47").Value = 'Investor''s Business Daily'
$ws.Range("F47").Value = 'Yahoo_Search'

$ws.Range("A48").Value = 46001
$ws.Range("B48").Value = 46001.77490740741
$ws.Range("C48").Value = 'Ethereum Has Bottomed, Says BitMine Chairman Tom Lee—Here''s Why'
$ws.Range("D48").Value = 'https://finance.yahoo.com/m/bcc9efc5-87f0-3f66-b601-094b802d156b/ethereum-has-bottomed%2C-says.html'
$ws.Range("E48").Value = 'decrypt'
$ws.Range("F48").Value = 'Yahoo_Search'

$ws.Range("A49").Value = 46001
$ws.Range("B49").Value = 46001.77430555555
$ws.Range("C49").Value = 'Three new Bachelors of the Politecnico di Milano, delivered entirely in English, in three cities'
$ws.Range("D49").Value = 'https://finance.yahoo.com/news/three-bachelors-politecnico-di-milano-183500230.html'
$ws.Range("E49").Value = 'GlobeNewswire'
$ws.Range("F49").Value = 'Yahoo_Search'

$ws.Range("A50").Value = 46001
$ws.Range("B50").Value = 46001.77430555555
$ws.Range("C50").Value = 'Mimeo Ltd. Acquires KnowledgePoint Print Services Ltd, Strengthening Its Leadership Position in Training Materials and Learning Fulfillment'
$ws.Range("D50").Value = 'https://finance.yahoo.com/news/mimeo-ltd-acquires-knowledgepoint-print-183500341.html'
$ws.Range("E50").Value = 'PR Newswire'
$ws.Range("F50").Value = 'Yahoo_Search'

# Apply date/datetime number formats to the newly added rows so they match
# the existing style used for columns A (date) and B (datetime).
$ws.Range("A42:A50").NumberFormat = "YYYY-MM-DD"
$ws.Range("B42:B50").NumberFormat = "YYYY-MM-DD HH:MM:SS"
